$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated covariates data (age, sex, edu) for GroupName1_covariates
$data = @(
    @(2, 57, "Male", 6),
    @(3, 60, "Male", 6),
    @(4, 61, "Female", 18),
    @(5, 60, "Female", 13),
    @(6, 67, "Female", 5),
    @(7, 76, "Male", 10),
    @(8, 78, "Female", 6),
    @(9, 73, "Male", 20),
    @(10, 59, "Male", 5),
    @(11, 61, "Female", 8),
    @(12, 52, "Male", 5),
    @(13, 53, "Male", 15),
    @(14, 51, "Female", 10),
    @(15, 76, "Male", 8),
    @(16, 50, "Male", 5),
    @(17, 61, "Female", 5),
    @(18, 68, "Male", 8),
    @(19, 57, "Female", 6),
    @(20, 74, "Female", 5),
    @(21, 73, "Female", 6),
    @(22, 63, "Male", 7),
    @(23, 71, "Female", 8),
    @(24, 71, "Female", 10),
    @(25, 61, "Female", 20),
    @(26, 63, "Male", 9),
    @(27, 69, "Female", 16),
    @(28, 78, "Male", 15),
    @(29, 75, "Male", 7),
    @(30, 66, "Male", 7),
    @(31, 65, "Female", 6),
    @(32, 59, "Female", 6),
    @(33, 60, "Male", 5),
    @(34, 63, "Female", 7),
    @(35, 58, "Male", 9),
    @(36, 78, "Female", 13),
    @(37, 77, "Male", 9),
    @(38, 72, "Male", 15),
    @(39, 63, "Male", 7),
    @(40, 56, "Female", 5),
    @(41, 59, "Male", 16),
    @(42, 52, "Male", 5),
    @(43, 79, "Male", 7),
    @(44, 60, "Female", 18),
    @(45, 57, "Female", 11),
    @(46, 54, "Female", 7),
    @(47, 72, "Male", 9),
    @(48, 70, "Male", 8),
    @(49, 54, "Male", 7),
    @(50, 56, "Male", 7),
    @(51, 60, "Female", 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
